# Update the "想去人数" (F column) counts on the "展览" and "全部类型"
# sheets to reflect newly generated output (gh-pages regeneration).

$wb = $excel.ActiveWorkbook

# Row -> new F-column value mapping (same updates apply to both sheets).
$updates = @{
    7  = 98
    10 = 130
    11 = 4532
    12 = 6799
    18 = 4123
    25 = 166
    35 = 81
    40 = 65
    42 = 12
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
